$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I6").Value = 611.875
$ws.Range("K6").Value = 1835.625
$ws.Range("L6").Value = 2908.8
$ws.Range("N6").Value = -3132.8
$ws.Range("M6").Value = -1723.625
$ws.Range("H6").Value = 749.46155
$ws.Range("J6").Value = 969.6
$ws.Range("K74").Value = 3689.6667
$ws.Range("M74").Value = -2753.6667
$ws.Range("H74").Value = 3267.25
$ws.Range("I74").Value = 3689.6667
$ws.Range("I77").Value = 3689.6667
$ws.Range("H77").Value = 3267.25
$ws.Range("M77").Value = -13768.3335
$ws.Range("K77").Value = 18448.3335
$ws.Range("H116").Value = 3301
$ws.Range("K116").Value = 3333.3333
$ws.Range("I116").Value = 3333.3333
$ws.Range("M116").Value = 108.6667000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 897.0909
$ws.Range("H2").Value = 897.0909
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -784.0909
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 897.0909
$ws.Range("N2").Value = $null
$ws.Range("L44").Value = 12441.353
$ws.Range("J44").Value = 12441.353
$ws.Range("H44").Value = 12441.353
$ws.Range("N44").Value = -13417.353
$ws.Range("K102").Value = 1570
$ws.Range("N102").Value = -11043.8
$ws.Range("M102").Value = 52
$ws.Range("L102").Value = 7799.8
$ws.Range("H102").Value = 3127.45
$ws.Range("I102").Value = 1570
$ws.Range("J102").Value = 7799.8
$ws.Range("H110").Value = 1389
$ws.Range("K110").Value = 1307.25
$ws.Range("M110").Value = 737.75
$ws.Range("I110").Value = 1307.25
$ws.Range("N110").Value = -5588
$ws.Range("L110").Value = 1498
$ws.Range("J110").Value = 1498
$ws.Range("H116").Value = 897.0909
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("K116").Value = 897.0909
$ws.Range("I116").Value = 897.0909
$ws.Range("N116").Value = $null
$ws.Range("M116").Value = 1396.9091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("I3").Value = 897.0909
$ws.Range("H3").Value = 897.0909
$ws.Range("K3").Value = 897.0909
$ws.Range("N3").Value = $null
$ws.Range("M3").Value = -783.0909
$ws.Range("K99").Value = 1602.2307
$ws.Range("H99").Value = 2101.2222
$ws.Range("I99").Value = 1602.2307
$ws.Range("M99").Value = -104.2307000000001
$ws.Range("H134").Value = 2139.4
$ws.Range("K134").Value = 2697.9999
$ws.Range("I134").Value = 899.3333
$ws.Range("M134").Value = -162.9998999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J22").Value = 2799
$ws.Range("H22").Value = 1383.9
$ws.Range("I22").Value = 777.4286
$ws.Range("M22").Value = -427.4286
$ws.Range("K22").Value = 777.4286
$ws.Range("N22").Value = -3499
$ws.Range("L22").Value = 2799
$ws.Range("M31").Value = -1688.08
$ws.Range("I31").Value = 1983.08
$ws.Range("H31").Value = 4685.421
$ws.Range("K31").Value = 1983.08
$ws.Range("H34").Value = 4685.421
$ws.Range("M34").Value = -1781.08
$ws.Range("I34").Value = 1983.08
$ws.Range("K34").Value = 1983.08
$ws.Range("J88").Value = 13285
$ws.Range("N88").Value = -14097
$ws.Range("H88").Value = 13285
$ws.Range("L88").Value = 13285
$ws.Range("J91").Value = 13285
$ws.Range("N91").Value = -16093
$ws.Range("L91").Value = 13285
$ws.Range("H91").Value = 13285
$ws.Range("H134").Value = 2185.7307
$ws.Range("K134").Value = 6601.875
$ws.Range("I134").Value = 2200.625
$ws.Range("M134").Value = -4066.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3931.1875
$ws.Range("L39").Value = 13343.3568
$ws.Range("J39").Value = 4447.7856
$ws.Range("N39").Value = -13931.3568
$ws.Range("N86").Value = -3869
$ws.Range("K86").Value = 949.1999999999999
$ws.Range("J86").Value = 499
$ws.Range("L86").Value = 1497
$ws.Range("H86").Value = 327.8125
$ws.Range("M86").Value = 236.8000000000001
$ws.Range("I86").Value = 316.4
$ws.Range("J89").Value = 499
$ws.Range("H89").Value = 327.8125
$ws.Range("L89").Value = 4491
$ws.Range("K89").Value = 2847.6
$ws.Range("I89").Value = 316.4
$ws.Range("M89").Value = 3080.4
$ws.Range("N89").Value = -16347
$ws.Range("H107").Value = 889.2857
$ws.Range("J107").Value = 926.9231
$ws.Range("N107").Value = -6620.7693
$ws.Range("L107").Value = 2780.7693
$ws.Range("I109").Value = 0
$ws.Range("M109").Value = $null
$ws.Range("K109").Value = 0
$ws.Range("H109").Value = 0
$ws.Range("H121").Value = 998
$ws.Range("J121").Value = 997.5
$ws.Range("L121").Value = 2992.5
$ws.Range("N121").Value = -5612.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K132").Value = 10123.5
$ws.Range("I132").Value = 3374.5
$ws.Range("M132").Value = -7593.5
$ws.Range("H132").Value = 3699.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I2").Value = 49999.5
$ws.Range("H2").Value = 79509
$ws.Range("M2").Value = -49887.5
$ws.Range("K2").Value = 49999.5
$ws.Range("I7").Value = 7646.5713
$ws.Range("M7").Value = -7534.5713
$ws.Range("K7").Value = 7646.5713
$ws.Range("H7").Value = 8038.8
$ws.Range("J22").Value = 302
$ws.Range("H22").Value = 302
$ws.Range("N22").Value = -892
$ws.Range("L22").Value = 302
$ws.Range("N27").Value = -516
$ws.Range("H27").Value = 302
$ws.Range("J27").Value = 302
$ws.Range("L27").Value = 302
$ws.Range("J40").Value = 3002.5
$ws.Range("K40").Value = 4043.889
$ws.Range("N40").Value = -3274.5
$ws.Range("H40").Value = 3854.5454
$ws.Range("I40").Value = 4043.889
$ws.Range("M40").Value = -3907.889
$ws.Range("L40").Value = 3002.5
$ws.Range("J122").Value = 3003.1667
$ws.Range("I122").Value = 3002
$ws.Range("L122").Value = 9009.500100000001
$ws.Range("N122").Value = -13909.5001
$ws.Range("K122").Value = 9006
$ws.Range("H122").Value = 3003
$ws.Range("M122").Value = -6556
$ws.Range("H126").Value = 8038.8
$ws.Range("I126").Value = 7646.5713
$ws.Range("K126").Value = 22939.7139
$ws.Range("M126").Value = -20469.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K62").Value = 6408.6665
$ws.Range("N62").Value = -11647.8
$ws.Range("J62").Value = 10399.8
$ws.Range("L62").Value = 10399.8
$ws.Range("H62").Value = 8222.817999999999
$ws.Range("M62").Value = -5784.6665
$ws.Range("I62").Value = 6408.6665
$ws.Range("K65").Value = 32043.3325
$ws.Range("J65").Value = 10399.8
$ws.Range("H65").Value = 8222.817999999999
$ws.Range("I65").Value = 6408.6665
$ws.Range("M65").Value = -28923.3325
$ws.Range("L65").Value = 51999
$ws.Range("N65").Value = -58239
$ws.Range("I107").Value = 3599
$ws.Range("H107").Value = 2849.25
$ws.Range("M107").Value = -8877
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 10797
$ws.Range("N107").Value = -5640
$ws.Range("L107").Value = 1800
$ws.Range("J122").Value = 7335
$ws.Range("I122").Value = 1485.6364
$ws.Range("L122").Value = 22005
$ws.Range("N122").Value = -26905
$ws.Range("K122").Value = 4456.9092
$ws.Range("H122").Value = 2739.0715
$ws.Range("M122").Value = -2006.9092

Write-Output "Applied 189 cell updates across 8 sheets"